# Apply the daily cryptos price/volume update (GitHub Actions data refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose text value would otherwise be re-interpreted by Excel as a
# number (e.g. "1.00" -> 1) are written with a leading apostrophe so they
# stay plain text, matching the original inline-string cells.

$ws.Range("D2").Value = '37.441.65'
$ws.Range("E2").Value = '  -0.95%  '

$ws.Range("D3").Value = '2.051.99'
$ws.Range("E3").Value = '  -1.67%  '

$ws.Range("D4").Value = '''1.00'
$ws.Range("E4").Value = '  +0.01%  '

$ws.Range("D5").Value = '''228.94'
$ws.Range("E5").Value = '  -2.07%  '

$ws.Range("D6").Value = '''0.613'
$ws.Range("E6").Value = '  -1.72%  '

$ws.Range("E7").Value = '  +0.05%  '

$ws.Range("D8").Value = '''56.35'
$ws.Range("E8").Value = '  -3.15%  '

$ws.Range("D9").Value = '''0.387'
$ws.Range("E9").Value = '  -1.69%  '

$ws.Range("E10").Value = '  +2.95%  '

$ws.Range("E12").Value = '  -1.74%  '

$ws.Range("D13").Value = '''14.56'
$ws.Range("E13").Value = '  -4.83%  '

$ws.Range("D14").Value = '''20.64'
$ws.Range("E14").Value = '  -2.52%  '

$ws.Range("D15").Value = '''0.755'
$ws.Range("E15").Value = '  -3.00%  '

$ws.Range("D16").Value = '''5.25'
$ws.Range("E16").Value = '  -2.10%  '

$ws.Range("D17").Value = '2.039.34'
$ws.Range("E17").Value = '  -2.35%  '

$ws.Range("D18").Value = '37.335.31'
$ws.Range("E18").Value = '  -1.07%  '

$ws.Range("E19").Value = '  -0.85%  '

$ws.Range("D20").Value = '''69.81'
$ws.Range("E20").Value = '  -1.70%  '

$ws.Range("D21").Value = '0.0₃0847'
$ws.Range("E21").Value = '  +1.42%  '

$ws.Range("D22").Value = '''225.81'
$ws.Range("E22").Value = '  -1.72%  '

$ws.Range("D23").Value = '''1.00'
$ws.Range("E23").Value = '  +0.08%  '

$ws.Range("E24").Value = '  -0.64%  '

$ws.Range("D25").Value = '''2.29'
$ws.Range("E25").Value = '  -4.43%  '

$ws.Range("D26").Value = '''9.50'

$ws.Range("D27").Value = '''168.14'
$ws.Range("E27").Value = '  -1.96%  '

$ws.Range("D28").Value = '''0.128'
$ws.Range("E28").Value = '  -4.65%  '

$ws.Range("E29").Value = '  -2.02%  '

$ws.Range("D30").Value = '''18.93'
$ws.Range("E30").Value = '  -2.83%  '

$ws.Range("D32").Value = '''4.54'
$ws.Range("E32").Value = '  -3.31%  '

$ws.Range("E33").Value = '  -3.26%  '

$ws.Range("D34").Value = '''4.53'
$ws.Range("E34").Value = '  -2.29%  '

$ws.Range("D35").Value = '''2.40'
$ws.Range("E35").Value = '  -3.95%  '

$ws.Range("E36").Value = '  +0.34%  '

$ws.Range("E37").Value = '  +0.20%  '

$ws.Range("D38").Value = '''3.20'
$ws.Range("E38").Value = '  -4.04%  '

$ws.Range("D39").Value = '''5.43'
$ws.Range("E39").Value = '  +0.87%  '

$ws.Range("E40").Value = '  -6.68%  '

$ws.Range("D41").Value = '1.495.01'
$ws.Range("E41").Value = '  +3.13%  '

$ws.Range("B42").Value = 'InjectiveProtocol'
$ws.Range("C42").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D42").Value = '''16.76'
$ws.Range("E42").Value = '  -0.17%  '

$ws.Range("B43").Value = 'Aave'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D43").Value = '''96.04'
$ws.Range("E43").Value = '  -5.07%  '

$ws.Range("B44").Value = 'HuobiToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D44").Value = '''2.85'
$ws.Range("E44").Value = '  -2.05%  '

$ws.Range("D45").Value = '''0.0935'
$ws.Range("E45").Value = '  -3.56%  '

$ws.Range("D46").Value = '''1.15'
$ws.Range("E46").Value = '  -3.96%  '

$ws.Range("D47").Value = '''1.02'
$ws.Range("E47").Value = '  -4.07%  '

$ws.Range("D48").Value = '''7.22'
$ws.Range("E48").Value = '  +0.18%  '

$ws.Range("D49").Value = '''2.93'
$ws.Range("E49").Value = '  -1.02%  '

$ws.Range("D50").Value = '''3.83'
$ws.Range("E50").Value = '  -7.23%  '

$ws.Range("D51").Value = '2.239.32'
$ws.Range("E51").Value = '  -1.71%  '
